{"js": "// Apply the two-digit multiplication worksheet update:\n// the date line and each \"AxB=\" cell get replaced with new values,\n// matching the published diff 1:1 (each old value is unique in the doc).\nconst replacements = [\n  [\"2024-09-08 Sunday\", \"2024-09-09 Monday\"],\n  [\"25\u00d747=\", \"50\u00d790=\"],\n  [\"41\u00d747=\", \"47\u00d718=\"],\n  [\"23\u00d738=\", \"84\u00d734=\"],\n  [\"86\u00d796=\", \"76\u00d780=\"],\n  [\"71\u00d720=\", \"84\u00d768=\"],\n  [\"65\u00d788=\", \"60\u00d739=\"],\n  [\"89\u00d782=\", \"39\u00d776=\"],\n  [\"79\u00d774=\", \"60\u00d787=\"],\n  [\"59\u00d766=\", \"53\u00d779=\"],\n  [\"31\u00d749=\", \"24\u00d758=\"],\n  [\"44\u00d746=\", \"19\u00d750=\"],\n  [\"36\u00d758=\", \"80\u00d764=\"],\n  [\"87\u00d793=\", \"17\u00d713=\"],\n  [\"56\u00d768=\", \"51\u00d794=\"],\n  [\"68\u00d712=\", \"74\u00d747=\"],\n  [\"31\u00d777=\", \"12\u00d765=\"],\n  [\"31\u00d763=\", \"57\u00d774=\"],\n  [\"90\u00d735=\", \"51\u00d771=\"],\n  [\"28\u00d727=\", \"87\u00d746=\"],\n  [\"20\u00d774=\", \"21\u00d798=\"],\n  [\"66\u00d788=\", \"82\u00d726=\"],\n  [\"52\u00d752=\", \"53\u00d774=\"],\n  [\"44\u00d717=\", \"62\u00d741=\"],\n  [\"88\u00d732=\", \"76\u00d745=\"],\n  [\"54\u00d780=\", \"94\u00d723=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, 'Replace');\n  }\n}\n\nawait context.sync();", "ps1": "# Update master to output generated at c986bee\n# Replaces the worksheet date line and each \"AxB=\" prompt with its new\n# value; every old value is unique in the document, so a plain\n# Find/Replace (no wildcards) on $d.Content relocates the right cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-09-08 Sunday\", \"2024-09-09 Monday\"),\n    @(\"25\u00d747=\", \"50\u00d790=\"),\n    @(\"41\u00d747=\", \"47\u00d718=\"),\n    @(\"23\u00d738=\", \"84\u00d734=\"),\n    @(\"86\u00d796=\", \"76\u00d780=\"),\n    @(\"71\u00d720=\", \"84\u00d768=\"),\n    @(\"65\u00d788=\", \"60\u00d739=\"),\n    @(\"89\u00d782=\", \"39\u00d776=\"),\n    @(\"79\u00d774=\", \"60\u00d787=\"),\n    @(\"59\u00d766=\", \"53\u00d779=\"),\n    @(\"31\u00d749=\", \"24\u00d758=\"),\n    @(\"44\u00d746=\", \"19\u00d750=\"),\n    @(\"36\u00d758=\", \"80\u00d764=\"),\n    @(\"87\u00d793=\", \"17\u00d713=\"),\n    @(\"56\u00d768=\", \"51\u00d794=\"),\n    @(\"68\u00d712=\", \"74\u00d747=\"),\n    @(\"31\u00d777=\", \"12\u00d765=\"),\n    @(\"31\u00d763=\", \"57\u00d774=\"),\n    @(\"90\u00d735=\", \"51\u00d771=\"),\n    @(\"28\u00d727=\", \"87\u00d746=\"),\n    @(\"20\u00d774=\", \"21\u00d798=\"),\n    @(\"66\u00d788=\", \"82\u00d726=\"),\n    @(\"52\u00d752=\", \"53\u00d774=\"),\n    @(\"44\u00d717=\", \"62\u00d741=\"),\n    @(\"88\u00d732=\", \"76\u00d745=\"),\n    @(\"54\u00d780=\", \"94\u00d723=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
